$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Time Type" column in M with header and values
$ws.Range("M1").Value = "Time Type"
$ws.Range("M2").Value = "Part Time"
$ws.Range("M3").Value = "Full Time"
$ws.Range("M4").Value = "Full Time"
$ws.Range("M5").Value = "Part Time"
$ws.Range("M6").Value = "Full Time"
$ws.Range("M7").Value = "Part Time"
$ws.Range("M8").Value = "Full Time"

# Update the selection to match the final state
$ws.Range("M7").Select()
